$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.697.83"
$ws.Range("E2").Value = "  -0.29%  "
$ws.Range("D3").Value = "2.396.88"
$ws.Range("E3").Value = "  -0.76%  "
$ws.Range("E4").Value = "  +0.41%  "
$ws.Range("D5").Value = "'563.09"
$ws.Range("E5").Value = "  -1.03%  "
$ws.Range("D6").Value = "'141.42"
$ws.Range("E6").Value = "  +1.75%  "
$ws.Range("E7").Value = "  -0.38%  "
$ws.Range("E8").Value = "  +2.37%  "
$ws.Range("D9").Value = "2.403.42"
$ws.Range("E9").Value = "  +0.16%  "
$ws.Range("E10").Value = "  +1.11%  "
$ws.Range("E11").Value = "  -0.21%  "
$ws.Range("D12").Value = "'5.19"
$ws.Range("E12").Value = "  +2.42%  "
$ws.Range("E13").Value = "  +2.24%  "
$ws.Range("D14").Value = "'26.30"
$ws.Range("E14").Value = "  +0.96%  "
$ws.Range("E15").Value = "  -1.03%  "
$ws.Range("D16").Value = "2.801.78"
$ws.Range("E16").Value = "  -2.05%  "
$ws.Range("D17").Value = "60.589.75"
$ws.Range("E17").Value = "  -0.40%  "
$ws.Range("D18").Value = "2.410.34"
$ws.Range("E18").Value = "  +0.64%  "
$ws.Range("D19").Value = "'8.07"
$ws.Range("E19").Value = "  +4.17%  "
$ws.Range("D20").Value = "'10.68"
$ws.Range("E20").Value = "  +0.44%  "
$ws.Range("D21").Value = "'324.11"
$ws.Range("E21").Value = "  +0.44%  "
$ws.Range("E22").Value = "  +1.45%  "
$ws.Range("D23").Value = "'6.05"
$ws.Range("E23").Value = "  -0.44%  "
$ws.Range("E24").Value = "  -0.20%  "
$ws.Range("D25").Value = "'1.87"
$ws.Range("E25").Value = "  +2.13%  "
$ws.Range("D26").Value = "'64.93"
$ws.Range("E26").Value = "  +0.11%  "
$ws.Range("D27").Value = "'582.09"
$ws.Range("E27").Value = "  +0.52%  "
$ws.Range("D28").Value = "'8.14"
$ws.Range("E28").Value = "  -2.13%  "
$ws.Range("D30").Value = "0.0₃0936"
$ws.Range("E30").Value = "  +1.15%  "
$ws.Range("E31").Value = "  +2.31%  "
$ws.Range("E32").Value = "  +0.33%  "
$ws.Range("E33").Value = "  -1.15%  "
$ws.Range("E34").Value = "  +0.34%  "
$ws.Range("E35").Value = "  -0.67%  "
$ws.Range("E36").Value = "  +4.69%  "
$ws.Range("B37").Value = "PolygonEcosystemToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D37").Value = "'0.371"
$ws.Range("E37").Value = "  +1.16%  "
$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").Value = "'4.60"
$ws.Range("E38").Value = "  +0.39%  "
$ws.Range("B39").Value = "Monero"
$ws.Range("C39").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D39").Value = "'151.15"
$ws.Range("E39").Value = "  +0.09%  "
$ws.Range("D40").Value = "'18.29"
$ws.Range("E40").Value = "  +0.43%  "
$ws.Range("D41").Value = "'5.16"
$ws.Range("E41").Value = "  +0.73%  "
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("E43").Value = "  +7.37%  "
$ws.Range("D44").Value = "'1.68"
$ws.Range("E44").Value = "  +0.89%  "
$ws.Range("D45").Value = "'41.60"
$ws.Range("E45").Value = "  +1.21%  "
$ws.Range("E46").Value = "  +5.46%  "
$ws.Range("D47").Value = "'141.12"
$ws.Range("E47").Value = "  -0.82%  "
$ws.Range("E48").Value = "  +0.46%  "
$ws.Range("E49").Value = "  +0.84%  "
$ws.Range("D51").Value = "'19.37"
$ws.Range("E51").Value = "  +0.10%  "
